$d = $word.ActiveDocument

$replacements = @(
    @{old="397×4=1588"; new="870×6=5220"},
    @{old="275×2=550"; new="237×6=1422"},
    @{old="361×4=1444"; new="356×8=2848"},
    @{old="922×9=8298"; new="995×2=1990"},
    @{old="316×7=2212"; new="534×5=2670"},
    @{old="735×6=4410"; new="329×7=2303"},
    @{old="583×8=4664"; new="196×5=980"},
    @{old="516×2=1032"; new="949×9=8541"},
    @{old="915×3=2745"; new="921×4=3684"},
    @{old="271×5=1355"; new="459×2=918"},
    @{old="403×2=806"; new="242×9=2178"},
    @{old="905×3=2715"; new="749×2=1498"},
    @{old="197×3=591"; new="692×7=4844"},
    @{old="735×4=2940"; new="707×8=5656"},
    @{old="516×6=3096"; new="315×8=2520"},
    @{old="780×9=7020"; new="515×6=3090"},
    @{old="182×2=364"; new="398×3=1194"},
    @{old="356×2=712"; new="325×7=2275"},
    @{old="424×4=1696"; new="465×2=930"},
    @{old="581×5=2905"; new="353×8=2824"},
    @{old="840×5=4200"; new="647×2=1294"},
    @{old="235×6=1410"; new="963×2=1926"},
    @{old="593×5=2965"; new="422×4=1688"},
    @{old="138×8=1104"; new="645×9=5805"},
    @{old="615×9=5535"; new="362×6=2172"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
